$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "css" -> "scss" in the task descriptions for rows 45-51 ---
# (row 52 keeps its original text "html done ja css pääosin" unchanged)
$ws.Range("C45").Value = "Html sivujen ja scss:n toteuttamista"
$ws.Range("C46").Value = "Landing page scss: footer ja content"
$ws.Range("C47").Value = "SCSS: responsiivisuutta"
$ws.Range("C48").Value = "SCSS: sivuja"
$ws.Range("C49").Value = "SCSS: responsiivisuutta paikka sivu"
$ws.Range("C50").Value = "SCSS: Responsiivisuus käyttäjä ja paikka sivu done "
$ws.Range("C51").Value = "SCSS: responsiivisuus käyttäjä asetus sivuja"

# --- Fill in the two new timesheet rows ---
$ws.Range("A53").Value = 44079
$ws.Range("B53").Value = 3
$ws.Range("C53").Value = "html:n ja scss:n siirto Reactiin"

$ws.Range("A54").Value = 44080
$ws.Range("B54").Value = 1
$ws.Range("C54").Value = "html:n ja scss:n siirto Reactiin + jakamista komponentteihin"

# Row 54's description is long, so it wraps (matching the style used by other
# multi-line rows such as C5/C8/C9) and the row is taller to fit it.
$ws.Range("C54").WrapText = $true
$ws.Rows.Item(54).RowHeight = 30

# --- Update the view: scrolled down one row further, new active cell ---
[void]$ws.Range("C51").Select()
